$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Replace the (too-short) second task description with the longer text ---
# Target cell: C3 ("The second task description" -> longer text). The text
# below uses a Unicode right single quotation mark (U+2019) in "that's" to
# match the source content exactly.
$ws.Range("C3").Value = "The second task description is way too long to be fully shown on the description part of the sticky, but it does not push the sprint id off the sticky because that’s in a separate table cell"

# --- Column width changes ---
# The host's ColumnWidth setter quantizes to whole pixels (1/6-character
# steps for this font) and the saved OOXML `width` attribute is that
# quantized character width plus a fixed 5-pixel (5/6 character) padding
# term, i.e. saved_width = round(ColumnWidth * 6)/6 + 5/6. The values below
# are chosen so that, after that transform, the saved width lands as close
# as possible to the widths from the target file:
#   col A       -> 8.36734693877551
#   col B       -> 14.3112244897959
#   col C       -> 23.0816326530612
#   cols D..AMK -> 8.36734693877551 (column 1025)
$ws.Columns.Item(1).ColumnWidth = 7.5
$ws.Columns.Item(2).ColumnWidth = 13.5
$ws.Columns.Item(3).ColumnWidth = 22.166666666666668
$ws.Range("D1:AMK1").EntireColumn.ColumnWidth = 7.5
